$wb = $excel.ActiveWorkbook

# Updates to apply, keyed by sheet name -> cell -> new value
$updates = @{
    "展览"     = @{ "F4" = 7935; "F5" = 5789; "F6" = 486; "F7" = 82; "F10" = 274; "F11" = 334 }
    "全部类型" = @{ "F4" = 7935; "F5" = 5789; "F6" = 486; "F7" = 82; "F10" = 274; "F14" = 334 }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellUpdates = $updates[$sheetName]
    foreach ($cellRef in $cellUpdates.Keys) {
        $ws.Range($cellRef).Value = $cellUpdates[$cellRef]
    }
}
